# Edit script: update cryptocurrency list data (rows 2-51) on sheet1.
#
# Applies the diff that refreshes the scraped cryptocurrency table
# (commit: "Updated cryptos list on Sun Jun 11 17:50:15 UTC 2023 with
# GitHub Actions"). The coin ranking shifted by one position (a new coin,
# Decentraland, was appended at the bottom of the list), and every
# Price / Volume(1h) cell was refreshed with newly scraped values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "25.921.26", "  -0.10%  "),
    @(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.754.48", "  -0.15%  "),
    @(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.000", "  +0.62%  "),
    @(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "236.27", "  -1.41%  "),
    @(6, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.000", "  +0.40%  "),
    @(7, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.5171", "  +2.85%  "),
    @(8, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.2694", "  +1.68%  "),
    @(9, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.06199", "  +0.12%  "),
    @(10, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.757.62", "  +0.09%  "),
    @(11, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.06993", "  +0.95%  "),
    @(12, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "15.45", "  -1.43%  "),
    @(13, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.6392", "  +7.54%  "),
    @(14, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "4.483", "  -0.27%  "),
    @(15, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "78.00", "  +0.30%  "),
    @(16, "BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "1.000", "  +0.86%  "),
    @(17, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "0.9984", "  +0.40%  "),
    @(18, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "25.946.15", "  -0.16%  "),
    @(19, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "11.65", "  -0.48%  "),
    @(20, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.000006691", "  -1.40%  "),
    @(21, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "1.981.06", "  +0.60%  "),
    @(22, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "4.075", "  -0.19%  "),
    @(23, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "8.374", "  +3.48%  "),
    @(24, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "5.180", "  +1.08%  "),
    @(25, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "136.35", "  -1.20%  "),
    @(26, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "1.486", "  -3.08%  "),
    @(27, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "15.16", "  +1.70%  "),
    @(28, "LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "1.830", "  -1.68%  "),
    @(29, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "103.21", "  +0.50%  "),
    @(30, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.08356", "  +3.39%  "),
    @(31, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "3.697", "  -2.13%  "),
    @(32, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "3.404", "  -2.08%  "),
    @(33, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.04393", "  -2.12%  "),
    @(34, "HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "2.643", "  +0.80%  "),
    @(35, "ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "0.9960", "  -0.42%  "),
    @(36, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "0.6061", "  -0.18%  "),
    @(37, "MXToken", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", "2.723", "  +0.92%  "),
    @(38, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.01564", "  +2.39%  "),
    @(39, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "1.944", "  -0.95%  "),
    @(40, "PaxDollar", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp", "0.9999", "  +0.53%  "),
    @(41, "Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "102.42", "  -2.28%  "),
    @(42, "TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "0.3866", "  +0.46%  "),
    @(43, "TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "0.7487", "  +1.69%  "),
    @(44, "FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "4.921", "  -4.70%  "),
    @(45, "Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "0.05495", "  +5.41%  "),
    @(46, "Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "0.1109", "  -0.59%  "),
    @(47, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "6.070", "  +1.28%  "),
    @(48, "Elrond", "https://coinranking.com/coin/omwkOTglq+elrond-egld", "30.25", "  -0.02%  "),
    @(49, "Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "52.75", "  +0.37%  "),
    @(50, "USDD", "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd", "1.004", "  +0.66%  "),
    @(51, "Decentraland", "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana", "0.3416", "  -1.13%  ")
)

foreach ($row in $rows) {
    $r = $row[0]

    # Coin name
    $ws.Cells.Item($r, 2).Value = $row[1]

    # Link
    $ws.Cells.Item($r, 3).Value = $row[2]

    # Price - prefix with a leading apostrophe so Excel stores values such
    # as "1.000" or "25.921.26" as literal text (matching the original
    # inline-string cell content) instead of re-parsing them as numbers.
    $ws.Cells.Item($r, 4).Value = "'" + $row[3]

    # Volume(1h)
    $ws.Cells.Item($r, 5).Value = $row[4]
}
